{"js": "// Insert three new bullet points after the paragraph ending in\n// \"MUCH cleaner, but may have too much trouble with missing methods in\n// heap dump graph (can't do type checking for example...)\" and before the\n// \"Arrays\" paragraph, matching the target diff.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its distinctive text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"MUCH cleaner, but may have too much trouble\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph ('MUCH cleaner...')\");\n}\n\n// New paragraphs to add, in order, each with its own list indent level.\nconst newItems = [\n  { text: \"For unmapped paths, act as if they don\\u2019t exist\", level: 3 },\n  { text: \"I.e. File.exist() just returns false rather than throwing an exception\", level: 4 },\n  {\n    text:\n      \"Seemed necessary when getting examples to run, but unfortunately causes \" +\n      \"URLClassLoader.findClass() (and similar) to fail silently if you\\u2019ve \" +\n      \"forgotten to map a class path\",\n    level: 5,\n  },\n];\n\n// Insert each paragraph right after the previous one, preserving order.\nlet previous = anchor;\nfor (const item of newItems) {\n  const newParagraph = previous.insertParagraph(item.text, \"After\");\n  newParagraph.listItem.level = item.level;\n  previous = newParagraph;\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet points after the paragraph ending in\n# \"MUCH cleaner, but may have too much trouble with missing methods in\n# heap dump graph (can't do type checking for example...)\" and before the\n# \"Arrays\" paragraph, matching the target diff.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*MUCH cleaner, but may have too much trouble*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find anchor paragraph ('MUCH cleaner...')\"\n}\n\n# New paragraphs to add, in order, each with its own list indent level\n# (ListLevelNumber is 1-based; level 4 -> ilvl=3, 5 -> ilvl=4, 6 -> ilvl=5).\n$newItems = @(\n    @{ Text = \"For unmapped paths, act as if they don\u2019t exist\"; Level = 4 },\n    @{ Text = \"I.e. File.exist() just returns false rather than throwing an exception\"; Level = 5 },\n    @{ Text = \"Seemed necessary when getting examples to run, but unfortunately causes URLClassLoader.findClass() (and similar) to fail silently if you\u2019ve forgotten to map a class path\"; Level = 6 }\n)\n\n$previous = $target\nforeach ($item in $newItems) {\n    $previous.Range.InsertParagraphAfter()\n    $newPara = $previous.Next()\n    $newPara.Range.Text = $item.Text\n    $newPara.Range.ListFormat.ListLevelNumber = $item.Level\n    $previous = $newPara\n}\n"}
